$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 8, pushing the existing rows 8-18 down to 10-20.
$ws.Rows("8:9").Insert()

# Apply the same date number format (used throughout column D) to the new date cells.
$ws.Range("D8:D9").NumberFormat = $ws.Range("D10").NumberFormat

# Row 8 (new): Castle Brite / Primera entry dated 2021-12-06
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44536
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100103
$ws.Range("H8").Value = "Frutos de hueso (carozo)"
$ws.Range("I8").Value = 100103003
$ws.Range("J8").Value = "Damasco"
$ws.Range("K8").Value = "Castle Brite"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 23000
$ws.Range("P8").Value = 22500
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1250
$ws.Range("T8").Value = 18

# Row 9 (new): Castle Brite / Segunda entry dated 2021-12-06
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44536
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103003
$ws.Range("J9").Value = "Damasco"
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("Q9").Value = "$/caja 18 kilos"
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 1028
$ws.Range("T9").Value = 18
